$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: bump the date in A1 by one day
$ws.Range("A1").Value = 45309

# Step 2: update the price list values
$ws.Range("D33").Value = 214.542
$ws.Range("D34").Value = 237.684
$ws.Range("D35").Value = 274.92
$ws.Range("D39").Value = 293.364
$ws.Range("D40").Value = 441.09
